$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 231 (open/low unchanged; high/close revised) ---
$ws.Cells.Item(231, 4).Value = 458.92   # D231 high
$ws.Cells.Item(231, 6).Value = 451.56   # F231 close

# --- Append new monthly rows 232-234 ---
$newRows = @(
    @{ Row = 232; Date = 45047.33333333334; Open = 444.76; High = 452.98; Low = 440.61; Close = 446.11; Volume = 0 },
    @{ Row = 233; Date = 45078.33333333334; Open = 446.11; High = 454.98; Low = 443.36; Close = 450.31; Volume = 0 },
    @{ Row = 234; Date = 45110.33333333334; Open = 450.31; High = 450.31; Low = 442.21; Close = 445.06; Volume = 0 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Column A carries the bordered/centered datetime style (s="2" on row 231) -
    # copy the source cell's formatting first, then overwrite with the new value.
    $ws.Cells.Item(231, 1).Copy($ws.Cells.Item($row, 1))
    $ws.Cells.Item($row, 1).Value = $r.Date

    $ws.Cells.Item($row, 2).Value = "FX_IDC:USDKZT"
    $ws.Cells.Item($row, 3).Value = $r.Open
    $ws.Cells.Item($row, 4).Value = $r.High
    $ws.Cells.Item($row, 5).Value = $r.Low
    $ws.Cells.Item($row, 6).Value = $r.Close
    $ws.Cells.Item($row, 7).Value = $r.Volume
}
